$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = 46004
$ws.Range("B26").Value = "四方坪站"
$ws.Range("C26").Value = 10330.68
$ws.Range("D26").Value = 8752.34
$ws.Range("E26").Value = 3388.73
$ws.Range("F26").Value = 428

$ws.Range("A27").Value = 46004
$ws.Range("B27").Value = "高岭站"
$ws.Range("C27").Value = 6216.39
$ws.Range("D27").Value = 5452.91
$ws.Range("E27").Value = 1624.58
$ws.Range("F27").Value = 213

$ws.Range("I27").Select()
